$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A93:D93").Copy($ws.Range("A94:D94"))
